# Update the pl_mw line-results sheet with the recalculated 380 kV case values
# (commit: "case with 380 kV done"). For each data row, columns B,D,E,F,G,H,K,M
# get new computed values; all other cells are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = @{ "B"=0.1423581687604241; "D"=0.02383153491422441; "E"=0.1418975379009098; "F"=1.001026938071632; "G"=0.8641787305874686; "H"=0.8706292266536622; "K"=0.5573233552341037; "M"=0.2719776128062676 }
    3 = @{ "B"=0.1329743783578579; "D"=0.02378178700569933; "E"=0.1339829978956288; "F"=0.9812616942354708; "G"=0.8438887487805147; "H"=0.866345540138056; "K"=0.4857134367956064; "M"=0.2428724109972933 }
    4 = @{ "B"=0.1272856977271033; "D"=0.02375095444888053; "E"=0.1292268806890888; "F"=0.9697991521135236; "G"=0.8320584939190212; "H"=0.8641968803054851; "K"=0.44168322027042; "M"=0.2250965958248798 }
    5 = @{ "B"=0.1249859700120908; "D"=0.02373832529260866; "E"=0.1273144455062223; "F"=0.9652967194621311; "G"=0.8273945898377093; "H"=0.8634420887254919; "K"=0.4237251842724845; "M"=0.2178763495970912 }
    6 = @{ "B"=0.1246052195724161; "D"=0.02373622455530189; "E"=0.126998432139267; "F"=0.964559259779179; "G"=0.8266296105515352; "H"=0.8633240442282215; "K"=0.4207423408903423; "M"=0.2166788459169453 }
    7 = @{ "B"=0.1272546079360524; "D"=0.02375078437946243; "E"=0.1292009851570626; "F"=0.9697377488022454; "G"=0.8319949602540078; "H"=0.8641862121696136; "K"=0.4414410940053131; "M"=0.2249991261519781 }
    8 = @{ "B"=0.1391075412599889; "D"=0.02381444585115844; "E"=0.1391469574943613; "F"=0.9940716974609387; "G"=0.8570519012207853; "H"=0.869052093908806; "K"=0.5326450310859627; "M"=0.2619222408201054 }
    9 = @{ "B"=0.1629274003651062; "D"=0.02393669614612648; "E"=0.159485942283716; "F"=1.047168177046728; "G"=0.911212834283333; "H"=0.8824292110718943; "K"=0.7110140711624808; "M"=0.3350996347987021 }
    10 = @{ "B"=0.1807771817683346; "D"=0.02402454695577205; "E"=0.1749592532295523; "F"=1.089509707326869; "G"=0.954131604163365; "H"=0.8946174944608742; "K"=0.8417902401787387; "M"=0.3893642377240241 }
    11 = @{ "B"=0.1889730715689808; "D"=0.02406401083388943; "E"=0.1821182961332042; "F"=1.109507169350934; "G"=0.9743494451220727; "H"=0.900679666193696; "K"=0.9012305361694359; "M"=0.4141665169787956 }
    12 = @{ "B"=0.1920874925238678; "D"=0.02407887703297718; "E"=0.1848468666574448; "F"=1.117186381708478; "G"=0.9821062257210258; "H"=0.9030500501382051; "K"=0.9237319720850792; "M"=0.4235757855534672 }
    13 = @{ "B"=0.1914162668818165; "D"=0.02407567888851858; "E"=0.184258432483972; "F"=1.115527772416314; "G"=0.9804311694843477; "H"=0.9025362147688156; "K"=0.9188862113695677; "M"=0.4215485612686933 }
    14 = @{ "B"=0.189229080693039; "D"=0.02406523547979589; "E"=0.1823424227852541; "F"=1.11013680184395; "G"=0.9749855760734079; "H"=0.9008731786744306; "K"=0.9030818897270194; "M"=0.4149402767332191 }
    15 = @{ "B"=0.1878907713998501; "D"=0.02405882827082451; "E"=0.1811711125636251; "F"=1.10684858271884; "G"=0.9716631380486263; "H"=0.8998642683387459; "K"=0.8934003323071806; "M"=0.4108947600679471 }
    16 = @{ "B"=0.1802430801987924; "D"=0.02402195731622569; "E"=0.1744938379235137; "F"=1.088217702940554; "G"=0.952824363959337; "H"=0.8942317622800147; "K"=0.8379046565183046; "M"=0.3877457356643248 }
    17 = @{ "B"=0.175570843956109; "D"=0.02399920573891379; "E"=0.1704285466559696; "F"=1.07697733059284; "G"=0.9414457759938273; "H"=0.8909092232381681; "K"=0.8038468900396651; "M"=0.3735747852437044 }
    18 = @{ "B"=0.1728906507695029; "D"=0.02398607291740262; "E"=0.1681015891638822; "F"=1.070581404216867; "G"=0.934966367555603; "H"=0.8890468913166103; "K"=0.7842530187005821; "M"=0.3654350241597228 }
    19 = @{ "B"=0.171984415207092; "D"=0.02398161853560765; "E"=0.167315650841708; "F"=1.068427722143568; "G"=0.9327837299691737; "H"=0.8884246927191271; "K"=0.7776180577006926; "M"=0.3626809186016402 }
    20 = @{ "B"=0.1760674718725141; "D"=0.02400163256238841; "E"=0.1708601330305797; "F"=1.078166715998051; "G"=0.9426502856854881; "H"=0.8912578700350764; "K"=0.8074728889931464; "M"=0.3750821666504791 }
    21 = @{ "B"=0.1898712175956661; "D"=0.02406830512262914; "E"=0.1829047211727683; "F"=1.111717361221409; "G"=0.9765823377291554; "H"=0.9013596207657599; "K"=0.9077242010125417; "M"=0.4168808208801238 }
    22 = @{ "B"=0.1989557556515393; "D"=0.02411142280149292; "E"=0.1908792700301589; "F"=1.134266391538389; "G"=0.9993464261690974; "H"=0.9083976774460325; "K"=0.973201578604403; "M"=0.4442989858766992 }
    23 = @{ "B"=0.1941014371629421; "D"=0.02408845377242841; "E"=0.1866135991754589; "F"=1.122174408596095; "G"=0.9871427319600627; "H"=0.9046013318221355; "K"=0.9382589886560311; "M"=0.4296560967639778 }
    24 = @{ "B"=0.1758429278896045; "D"=0.02400053555778214; "E"=0.1706649809855207; "F"=1.077628788343731; "G"=0.9421055328810723; "H"=0.8911000978117158; "K"=0.8058336172598501; "M"=0.3744006567916358 }
    25 = @{ "B"=0.1564220090277786; "D"=0.02390394288119246; "E"=0.1538921311309736; "F"=1.032222966839086; "G"=0.8960163930459544; "H"=0.8783974892914728; "K"=0.6628098569714496; "M"=0.3152174590851331 }
}

foreach ($row in $newValues.Keys) {
    foreach ($col in $newValues[$row].Keys) {
        $ws.Range("$col$row").Value = $newValues[$row][$col]
    }
}
